# "updating EDA Density Hierarchical"
# Update the numeric benchmark results and runtime/params annotations for
# the Agglo (row 3), BIRCH (row 4) and DBSCAN (row 5) algorithms, then
# move the sheet's selection/top-left to where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Agglo -----------------------------------------------------
$ws.Range("C3").Value = 81291.723
$ws.Range("D3").Value = 9.574
$ws.Range("E3").Value = 0.21
$ws.Range("F3").Value = 0.188
$ws.Range("G3").Value = 0.514
$ws.Range("H3").Value = "2.04 ± 0.06"
$ws.Range("I3").Value = "{'standard_scaler': False, 'pca__n_components': 0.8, 'agglomerative__affinity': 'euclidean', 'agglomerative__linkage': 'ward'}"

# --- Row 4: BIRCH -------------------------------------------------------
$ws.Range("C4").Value = 19257.202
$ws.Range("D4").Value = 9.005
$ws.Range("E4").Value = 0.279
$ws.Range("F4").Value = 0.27
$ws.Range("G4").Value = 0.522
$ws.Range("H4").Value = "0.85 ± 0.03"

# --- Row 5: DBSCAN -------------------------------------------------------
$ws.Range("C5").Value = 129112.436
$ws.Range("D5").Value = 11.153
$ws.Range("E5").Value = 0.497
$ws.Range("H5").Value = "2.82 ± 0.08"

# --- Selection / scroll position left by the author ---------------------
$ws.Range("H5").Select()
